$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell D values occasionally look like numbers (e.g. "1.00", "589.28").
# Force text formatting before assigning so Excel keeps the exact string
# (trailing zeros, no float rounding) instead of silently converting to a
# number; ClearFormats() afterwards drops the temporary text numberformat
# so the cell style matches the original (no explicit style applied).

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "63.152.07"
$c.ClearFormats()
$ws.Range("E2").Value = "  -0.03%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.577.03"
$c.ClearFormats()
$ws.Range("E3").Value = "  +0.28%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "589.28"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.73%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "144.80"
$c.ClearFormats()
$ws.Range("E6").Value = "  -1.77%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.590"
$c.ClearFormats()
$ws.Range("E8").Value = "  -1.86%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.107"
$c.ClearFormats()
$ws.Range("E9").Value = "  -1.91%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "5.59"
$c.ClearFormats()
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("E11").Value = "  -0.13%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.350"
$c.ClearFormats()
$ws.Range("E12").Value = "  -1.67%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "27.17"
$c.ClearFormats()
$ws.Range("E13").Value = "  -0.93%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.040.04"
$c.ClearFormats()
$ws.Range("E14").Value = "  +0.32%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "63.060.00"
$c.ClearFormats()
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("E16").Value = "  -1.48%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.577.41"
$c.ClearFormats()
$ws.Range("E17").Value = "  -0.24%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "11.09"
$c.ClearFormats()
$ws.Range("E18").Value = "  -2.35%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "341.51"
$c.ClearFormats()
$ws.Range("E19").Value = "  -0.66%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.34"
$c.ClearFormats()
$ws.Range("E20").Value = "  -1.77%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.66"
$c.ClearFormats()
$ws.Range("E21").Value = "  -3.33%  "
$ws.Range("E22").Value = "  +0.02%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.74"
$c.ClearFormats()
$ws.Range("E23").Value = "  +3.44%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "67.81"
$c.ClearFormats()
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("E25").Value = "  +8.11%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.62"
$c.ClearFormats()
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("E27").Value = "  -3.35%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.97"
$c.ClearFormats()
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("E29").Value = "  -0.12%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "8.27"
$c.ClearFormats()
$ws.Range("E30").Value = "  -2.65%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.94"
$c.ClearFormats()
$ws.Range("E31").Value = "  -2.63%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "466.37"
$c.ClearFormats()
$ws.Range("E32").Value = "  -0.19%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0₃0800"
$c.ClearFormats()
$ws.Range("E33").Value = "  -3.21%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.68"
$c.ClearFormats()
$ws.Range("E34").Value = "  +2.72%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "176.77"
$c.ClearFormats()
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  +0.05%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.397"
$c.ClearFormats()
$ws.Range("E37").Value = "  -3.09%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "18.87"
$c.ClearFormats()
$ws.Range("E38").Value = "  -1.90%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "4.58"
$c.ClearFormats()
$ws.Range("E39").Value = "  +0.75%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.06%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.70"
$c.ClearFormats()
$ws.Range("E41").Value = "  -3.09%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "158.50"
$c.ClearFormats()
$ws.Range("E42").Value = "  +4.33%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "40.00"
$c.ClearFormats()
$ws.Range("E43").Value = "  +0.98%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.70"
$c.ClearFormats()
$ws.Range("E44").Value = "  -2.92%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "21.38"
$c.ClearFormats()
$ws.Range("E45").Value = "  +1.91%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.633"
$c.ClearFormats()
$ws.Range("E46").Value = "  +3.29%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0538"
$c.ClearFormats()
$ws.Range("E47").Value = "  -1.59%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0963"
$c.ClearFormats()
$ws.Range("E48").Value = "  -1.70%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0237"
$c.ClearFormats()
$ws.Range("E49").Value = "  -1.13%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "18.14"
$c.ClearFormats()
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("E51").Value = "  +0.08%  "
